# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors: header cells AD1:AF1 styled like the other header cells, and
# AD:AF filled with the team's season record (85 wins, 77 losses, 0 ties)
# for every player row (2-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 44

# --- Headers (row 1): copy formatting from an existing header cell (AC1)
# so the new header cells pick up the same style (bold, centered, bordered)
# without registering a brand-new style record.
$ws.Cells.Item(1, 29).Copy()
$ws.Cells.Item(1, 30).PasteSpecial(-4122)
$ws.Cells.Item(1, 31).PasteSpecial(-4122)
$ws.Cells.Item(1, 32).PasteSpecial(-4122)

$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# --- Data rows (2-44): same season-record values for every player.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
